$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 145
$ws.Range("I4").Value = 114.6
$ws.Range("K4").Value = 114.6
$ws.Range("M4").Value = -0.5999999999999943
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("N112").ClearContents()
$ws.Range("H132").Value = 2214.7368
$ws.Range("I132").Value = 1880
$ws.Range("K132").Value = 5640
$ws.Range("M132").Value = -3110
$ws.Range("H138").Value = 3026.0334
$ws.Range("J138").Value = 3194.9092
$ws.Range("L138").Value = 9584.7276
$ws.Range("N138").Value = -19864.7276
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H122").Value = 1471.1428
$ws.Range("I122").Value = 1079.8
$ws.Range("K122").Value = 3239.4
$ws.Range("M122").Value = -789.3999999999996
$ws.Range("H132").Value = 1526.6666
$ws.Range("I132").Value = 1526.6666
$ws.Range("K132").Value = 4579.9998
$ws.Range("M132").Value = -2049.9998
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5399.5713
$ws.Range("I86").Value = 3499.75
$ws.Range("K86").Value = 3499.75
$ws.Range("M86").Value = -2376.75
$ws.Range("H89").Value = 5399.5713
$ws.Range("I89").Value = 3499.75
$ws.Range("K89").Value = 17498.75
$ws.Range("M89").Value = -11882.75
$ws.Range("H94").Value = 4916.5
$ws.Range("J94").Value = 4999.8
$ws.Range("L94").Value = 4999.8
$ws.Range("N94").Value = -5901.8
$ws.Range("H99").Value = 2416.6667
$ws.Range("I99").Value = 2625
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 2625
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -1127
$ws.Range("N99").Value = -4996
$ws.Range("H107").Value = 1882.6666
$ws.Range("I107").Value = 1837.5
$ws.Range("K107").Value = 1837.5
$ws.Range("M107").Value = 82.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3824
$ws.Range("I99").Value = 4098.6665
$ws.Range("K99").Value = 4098.6665
$ws.Range("M99").Value = -2600.6665
$ws.Range("H126").Value = 3824
$ws.Range("I126").Value = 4098.6665
$ws.Range("K126").Value = 12295.9995
$ws.Range("M126").Value = -9825.999500000002
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2000
$ws.Range("I70").Value = 2000
$ws.Range("K70").Value = 6000
$ws.Range("M70").Value = -5685
$ws.Range("H73").Value = 2000
$ws.Range("I73").Value = 2000
$ws.Range("K73").Value = 6000
$ws.Range("M73").Value = -4908
$ws.Range("H107").Value = 1150.5
$ws.Range("J107").Value = 1298
$ws.Range("L107").Value = 3894
$ws.Range("N107").Value = -7734
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 35184.832
$ws.Range("J20").Value = 35184.832
$ws.Range("L20").Value = 35184.832
$ws.Range("N20").Value = -35674.832
$ws.Range("H52").Value = 99999
$ws.Range("J52").Value = 99999
$ws.Range("L52").Value = 99999
$ws.Range("N52").Value = -100517
$ws.Range("H80").Value = 4250
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 3500
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -2502
$ws.Range("N80").Value = -6996
$ws.Range("H83").Value = 4250
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 17500
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -12508
$ws.Range("N83").Value = -34984
$ws.Range("H102").Value = 2225.818
$ws.Range("I102").Value = 2098.4
$ws.Range("K102").Value = 2098.4
$ws.Range("M102").Value = -476.4000000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 7334.6665
$ws.Range("I20").Value = 1002.5
$ws.Range("J20").Value = 19999
$ws.Range("K20").Value = 1002.5
$ws.Range("L20").Value = 19999
$ws.Range("M20").Value = -776.5
$ws.Range("N20").Value = -20451
$ws.Range("H22").Value = 4391.778
$ws.Range("J22").Value = 4503.7144
$ws.Range("L22").Value = 4503.7144
$ws.Range("N22").Value = -5093.7144
$ws.Range("H27").Value = 4391.778
$ws.Range("J27").Value = 4503.7144
$ws.Range("L27").Value = 4503.7144
$ws.Range("N27").Value = -4717.7144
$ws.Range("H42").Value = 10000
$ws.Range("I42").Value = 10000
$ws.Range("K42").Value = 10000
$ws.Range("M42").Value = -9437
$ws.Range("H43").Value = 18333
$ws.Range("I43").Value = 17499.5
$ws.Range("J43").Value = 20000
$ws.Range("K43").Value = 17499.5
$ws.Range("L43").Value = 20000
$ws.Range("M43").Value = -17306.5
$ws.Range("N43").Value = -20386
$ws.Range("H46").Value = 4666.6665
$ws.Range("J46").Value = 5625
$ws.Range("L46").Value = 5625
$ws.Range("N46").Value = -6001
$ws.Range("H49").Value = 10000
$ws.Range("I49").Value = 10000
$ws.Range("K49").Value = 10000
$ws.Range("M49").Value = -9853
$ws.Range("H55").Value = 727.4
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 727.4
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 727.4
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -1073.4
$ws.Range("H61").Value = 1094
$ws.Range("I61").Value = 1094
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1094
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -892
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 3000
$ws.Range("I68").Value = 3000
$ws.Range("K68").Value = 3000
$ws.Range("M68").Value = -2251
$ws.Range("H71").Value = 3000
$ws.Range("I71").Value = 3000
$ws.Range("K71").Value = 15000
$ws.Range("M71").Value = -11256
$ws.Range("H113").Value = 1094
$ws.Range("I113").Value = 1094
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1094
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1076
$ws.Range("N113").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 999
$ws.Range("J15").Value = 999
$ws.Range("L15").Value = 999
$ws.Range("N15").Value = -1575
$ws.Range("H96").Value = 1452.421
$ws.Range("I96").Value = 1359.7333
$ws.Range("J96").Value = 1800
$ws.Range("K96").Value = 1359.7333
$ws.Range("L96").Value = 1800
$ws.Range("M96").Value = 13.2666999999999
$ws.Range("N96").Value = -4546
$ws.Range("H126").Value = 1758.8
$ws.Range("I126").Value = 1222.7142
$ws.Range("K126").Value = 3668.1426
$ws.Range("M126").Value = -1198.1426
